# Enhance Siege Analytics "Full-Stack Development and Data Engineering" bullet list
# by inserting new bullet points describing the boundary estimation algorithm work.

$d = $word.ActiveDocument

# Locate the paragraph that introduces the Siege Analytics bullet list.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "Full-Stack Development and Data Engineering") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Could not find anchor paragraph 'Full-Stack Development and Data Engineering'"
}
else {
    $newTexts = @(
        "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
        "• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times",
        "• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis",
        "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
    )

    $cur = $target
    foreach ($t in $newTexts) {
        $cur.Range.InsertParagraphAfter()
        $cur = $cur.Next()
        $cur.Range.Text = $t
    }

    Write-Host "Inserted $($newTexts.Count) new bullet paragraphs after 'Full-Stack Development and Data Engineering'"
}
